# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (cloned from the "2021-Q4" sheet so it
#    inherits the same column layout / header styling) positioned right
#    before the "总计" (total) sheet, then overwrite its values.
# 2. Insert a new top data row into the "总计" sheet summarising 2022-Q1
#    (holding count = 4, holding value = 0.06 yi), pushing the existing
#    rows down and re-numbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet by cloning "2021-Q4", placed before
# the "总计" sheet so the tab order becomes: ... 2021-Q4, 2022-Q1, 总计
# ---------------------------------------------------------------------
$q4_2021 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$q4_2021.Copy($total)
$q1_2022 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1_2022.Name = "2022-Q1"

# Header row (row 1) text stays the same as 2021-Q4's header, so nothing
# to change there. Make sure we have exactly 4 data rows (rows 2-5); the
# cloned sheet currently has 2 data rows (rows 2-3), so copy the
# formatting of row 2 down into rows 4-5 before filling in values.
$q1_2022.Range("A2:H2").Copy()
$q1_2022.Range("A4:H5").PasteSpecial(-4122) | Out-Null
$q1_2022.Range("A1").Select() | Out-Null

# Fund code (B), fund scale / position columns (D:G) are numeric-looking
# text in this workbook (e.g. "501310", "94.80") - force text storage so
# leading/trailing zeros survive, matching the rest of the workbook.
$q1_2022.Range("B2:B5").NumberFormat = "@"
$q1_2022.Range("D2:G5").NumberFormat = "@"

$q1_2022.Range("A2").Value = 0
$q1_2022.Range("B2").Value = "501310"
$q1_2022.Range("C2").Value = "华宝标普沪港深中国增强价值指数（LOF）A"
$q1_2022.Range("D2").Value = "1.44"
$q1_2022.Range("E2").Value = "94.80"
$q1_2022.Range("F2").Value = "2.80"
$q1_2022.Range("G2").Value = "0.0403"
$q1_2022.Range("H2").Value = 10

$q1_2022.Range("A3").Value = 1
$q1_2022.Range("B3").Value = "001942"
$q1_2022.Range("C3").Value = "前海开源沪港深汇鑫灵活配置混合A"
$q1_2022.Range("D3").Value = "0.10"
$q1_2022.Range("E3").Value = "90.39"
$q1_2022.Range("F3").Value = "7.42"
$q1_2022.Range("G3").Value = "0.0074"
$q1_2022.Range("H3").Value = 6

$q1_2022.Range("A4").Value = 2
$q1_2022.Range("B4").Value = "001943"
$q1_2022.Range("C4").Value = "前海开源沪港深汇鑫灵活配置混合C"
$q1_2022.Range("D4").Value = "0.08"
$q1_2022.Range("E4").Value = "90.39"
$q1_2022.Range("F4").Value = "7.42"
$q1_2022.Range("G4").Value = "0.0059"
$q1_2022.Range("H4").Value = 6

$q1_2022.Range("A5").Value = 3
$q1_2022.Range("B5").Value = "007397"
$q1_2022.Range("C5").Value = "华宝标普沪港深中国增强价值指数（LOF）C"
$q1_2022.Range("D5").Value = "0.09"
$q1_2022.Range("E5").Value = "94.80"
$q1_2022.Range("F5").Value = "2.80"
$q1_2022.Range("G5").Value = "0.0025"
$q1_2022.Range("H5").Value = 10

# ---------------------------------------------------------------------
# Step 2: prepend a 2022-Q1 summary row to the "总计" sheet, pushing the
# existing rows (2021-Q4 .. 2020-Q4) down by one and re-indexing column A.
# (PasteSpecial(xlPasteAll) does not reliably carry the A-column "index"
# style onto brand-new cells in this host, so formats are copied down
# separately from the values, which are written explicitly below.)
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A2:D6").Copy() | Out-Null
$total.Range("A3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.06

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.05

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.16

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.13

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 4
$total.Range("D6").Value = 0.18

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 8
$total.Range("D7").Value = 0.3

$total.Range("A1").Select() | Out-Null
